# Update generated stats for 江西-漫展信息.xlsx (gh-pages output regeneration)
# Touches sheet "展览", sheet "演出" and sheet "全部类型" (sheet "本地生活" untouched).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: 展览 (Exhibitions)
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

$wsExpo.Range("F2").Value  = 622
$wsExpo.Range("F5").Value  = 4559
$wsExpo.Range("F7").Value  = 126

# Row 8: 新余·LD03盛夏动漫嘉年华 was cancelled.
$wsExpo.Range("C8").Value  = "新余·LD03盛夏动漫嘉年华（取消）"
$wsExpo.Range("F8").Value  = 125
$wsExpo.Range("G8").Value  = "不可售"

$wsExpo.Range("F9").Value  = 3077
$wsExpo.Range("F13").Value = 596
$wsExpo.Range("F14").Value = 512
$wsExpo.Range("F15").Value = 513
$wsExpo.Range("F16").Value = 361
$wsExpo.Range("F18").Value = 1759
$wsExpo.Range("F19").Value = 1305
$wsExpo.Range("F20").Value = 116
$wsExpo.Range("F21").Value = 1551
$wsExpo.Range("F25").Value = 525
$wsExpo.Range("F28").Value = 88
$wsExpo.Range("F31").Value = 3553
$wsExpo.Range("F32").Value = 743
$wsExpo.Range("F34").Value = 249
$wsExpo.Range("F36").Value = 1714

# ---------------------------------------------------------------------------
# Sheet: 演出 (Performances)
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")

$wsShow.Range("F3").Value = 38

# ---------------------------------------------------------------------------
# Sheet: 全部类型 (All types)
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("F2").Value  = 622
$wsAll.Range("F5").Value  = 4559
$wsAll.Range("F7").Value  = 126

# Row 8: 新余·LD03盛夏动漫嘉年华 was cancelled.
$wsAll.Range("C8").Value  = "新余·LD03盛夏动漫嘉年华（取消）"
$wsAll.Range("F8").Value  = 125
$wsAll.Range("G8").Value  = "不可售"

$wsAll.Range("F9").Value  = 3077
$wsAll.Range("F13").Value = 596
$wsAll.Range("F14").Value = 512
$wsAll.Range("F15").Value = 513
$wsAll.Range("F17").Value = 361
$wsAll.Range("F19").Value = 1759
$wsAll.Range("F20").Value = 1305
$wsAll.Range("F21").Value = 116
$wsAll.Range("F22").Value = 1551
$wsAll.Range("F26").Value = 525
$wsAll.Range("F29").Value = 88
$wsAll.Range("F32").Value = 3553
$wsAll.Range("F33").Value = 38
$wsAll.Range("F34").Value = 743
$wsAll.Range("F36").Value = 249
$wsAll.Range("F38").Value = 1715

Write-Host "Applied gh-pages data refresh (456a3b4)."
